$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.423038333333333
$ws.Range("H2").Value = 7.269115
$ws.Range("I2").Value = 0.002232822326725897
$ws.Range("J2").Value = 0.002232822326725897
$ws.Range("M2").Value = 1.050658666666667
$ws.Range("N2").Value = 3.151976
$ws.Range("O2").Value = 0.1514279810580986
$ws.Range("P2").Value = 0.1514279810580986
$ws.Range("Q2").Value = 2.545786224582223
$ws.Range("R2").Value = 22.91207602124
$ws.Range("S2").Value = 0.0003381117769975488
$ws.Range("T2").Value = 0.0003381117769975486
$ws.Range("G3").Value = 2.423038333333333
$ws.Range("H3").Value = 7.269115
$ws.Range("I3").Value = 0.002232822326725897
$ws.Range("J3").Value = 0.002232822326725897
$ws.Range("M3").Value = 0.8776213333333334
$ws.Range("N3").Value = 2.632864
$ws.Range("O3").Value = 0.126488678822602
$ws.Range("P3").Value = 0.126488678822602
$ws.Range("Q3").Value = 2.126510132817778
$ws.Range("R3").Value = 19.13859119536
$ws.Range("S3").Value = 0.0002824267461531668
$ws.Range("T3").Value = 0.0002824267461531667
$ws.Range("G4").Value = 2.423038333333333
$ws.Range("H4").Value = 7.269115
$ws.Range("I4").Value = 0.002232822326725897
$ws.Range("J4").Value = 0.002232822326725897
$ws.Range("M4").Value = 5.010059000000001
$ws.Range("N4").Value = 15.030177
$ws.Range("O4").Value = 0.7220833401192995
$ws.Range("P4").Value = 0.7220833401192994
$ws.Range("Q4").Value = 12.13956500926167
$ws.Range("R4").Value = 109.256085083355
$ws.Range("S4").Value = 0.001612283803575182
$ws.Range("T4").Value = 0.001612283803575181
$ws.Range("I5").Value = 0.00002132067390877311
$ws.Range("J5").Value = 0.0000213206739087731
$ws.Range("M5").Value = 1.050658666666667
$ws.Range("N5").Value = 3.151976
$ws.Range("O5").Value = 0.1514279810580986
$ws.Range("P5").Value = 0.1514279810580986
$ws.Range("Q5").Value = 0.02430908957066667
$ws.Range("R5").Value = 0.218781806136
$ws.Range("S5").Value = 0.000003228546604803591
$ws.Range("T5").Value = 0.00000322854660480359
$ws.Range("I6").Value = 0.00002132067390877311
$ws.Range("J6").Value = 0.0000213206739087731
$ws.Range("M6").Value = 0.8776213333333334
$ws.Range("N6").Value = 2.632864
$ws.Range("O6").Value = 0.126488678822602
$ws.Range("P6").Value = 0.126488678822602
$ws.Range("Q6").Value = 0.02030552478933334
$ws.Range("R6").Value = 0.182749723104
$ws.Range("S6").Value = 0.000002696823874328231
$ws.Range("T6").Value = 0.000002696823874328231
$ws.Range("I7").Value = 0.00002132067390877311
$ws.Range("J7").Value = 0.0000213206739087731
$ws.Range("M7").Value = 5.010059000000001
$ws.Range("N7").Value = 15.030177
$ws.Range("O7").Value = 0.7220833401192995
$ws.Range("P7").Value = 0.7220833401192994
$ws.Range("Q7").Value = 0.115917735083
$ws.Range("R7").Value = 1.043259615747
$ws.Range("S7").Value = 0.00001539530342964129
$ws.Range("T7").Value = 0.00001539530342964128
$ws.Range("G8").Value = 280.561096
$ws.Range("H8").Value = 841.683288
$ws.Range("I8").Value = 0.2585361818431078
$ws.Range("J8").Value = 0.2585361818431078
$ws.Range("M8").Value = 1.050658666666667
$ws.Range("N8").Value = 3.151976
$ws.Range("O8").Value = 0.1514279810580986
$ws.Range("P8").Value = 0.1514279810580986
$ws.Range("Q8").Value = 294.7739470418987
$ws.Range("R8").Value = 2652.965523377088
$ws.Range("S8").Value = 0.03914961204697127
$ws.Range("T8").Value = 0.03914961204697126
$ws.Range("G9").Value = 280.561096
$ws.Range("H9").Value = 841.683288
$ws.Range("I9").Value = 0.2585361818431078
$ws.Range("J9").Value = 0.2585361818431078
$ws.Range("M9").Value = 0.8776213333333334
$ws.Range("N9").Value = 2.632864
$ws.Range("O9").Value = 0.126488678822602
$ws.Range("P9").Value = 0.126488678822602
$ws.Range("Q9").Value = 246.2264031529813
$ws.Range("R9").Value = 2216.037628376832
$ws.Range("S9").Value = 0.03270190006917468
$ws.Range("T9").Value = 0.03270190006917468
$ws.Range("G10").Value = 280.561096
$ws.Range("H10").Value = 841.683288
$ws.Range("I10").Value = 0.2585361818431078
$ws.Range("J10").Value = 0.2585361818431078
$ws.Range("M10").Value = 5.010059000000001
$ws.Range("N10").Value = 15.030177
$ws.Range("O10").Value = 0.7220833401192995
$ws.Range("P10").Value = 0.7220833401192994
$ws.Range("Q10").Value = 1405.627644064664
$ws.Range("R10").Value = 12650.64879658198
$ws.Range("S10").Value = 0.1866846697269619
$ws.Range("T10").Value = 0.1866846697269619
$ws.Range("E11").Value = 2.0
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.028388
$ws.Range("H11").Value = 0.085164
$ws.Range("I11").Value = 0.00002615945416096516
$ws.Range("J11").Value = 0.00002615945416096516
$ws.Range("M11").Value = 1.050658666666667
$ws.Range("N11").Value = 3.151976
$ws.Range("O11").Value = 0.1514279810580986
$ws.Range("P11").Value = 0.1514279810580986
$ws.Range("Q11").Value = 0.02982609822933334
$ws.Range("R11").Value = 0.268434884064
$ws.Range("S11").Value = 0.000003961273329176831
$ws.Range("T11").Value = 0.00000396127332917683
$ws.Range("E12").Value = 2.0
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.028388
$ws.Range("H12").Value = 0.085164
$ws.Range("I12").Value = 0.00002615945416096516
$ws.Range("J12").Value = 0.00002615945416096516
$ws.Range("M12").Value = 0.8776213333333334
$ws.Range("N12").Value = 2.632864
$ws.Range("O12").Value = 0.126488678822602
$ws.Range("P12").Value = 0.126488678822602
$ws.Range("Q12").Value = 0.02491391441066667
$ws.Range("R12").Value = 0.224225229696
$ws.Range("S12").Value = 0.0000033088747955409
$ws.Range("T12").Value = 0.0000033088747955409
$ws.Range("E13").Value = 2.0
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.028388
$ws.Range("H13").Value = 0.085164
$ws.Range("I13").Value = 0.00002615945416096516
$ws.Range("J13").Value = 0.00002615945416096516
$ws.Range("M13").Value = 5.010059000000001
$ws.Range("N13").Value = 15.030177
$ws.Range("O13").Value = 0.7220833401192995
$ws.Range("P13").Value = 0.7220833401192994
$ws.Range("Q13").Value = 0.142225554892
$ws.Range("R13").Value = 1.280029994028
$ws.Range("S13").Value = 0.00001888930603624743
$ws.Range("T13").Value = 0.00001888930603624743
$ws.Range("G14").Value = 683.8555499999999
$ws.Range("H14").Value = 2051.56665
$ws.Range("I14").Value = 0.630170773317832
$ws.Range("J14").Value = 0.630170773317832
$ws.Range("M14").Value = 1.050658666666667
$ws.Range("N14").Value = 3.151976
$ws.Range("O14").Value = 0.1514279810580986
$ws.Range("P14").Value = 0.1514279810580986
$ws.Range("Q14").Value = 718.4987603556
$ws.Range("R14").Value = 6466.4888432004
$ws.Range("S14").Value = 0.09542548792534002
$ws.Range("T14").Value = 0.09542548792533999
$ws.Range("G15").Value = 683.8555499999999
$ws.Range("H15").Value = 2051.56665
$ws.Range("I15").Value = 0.630170773317832
$ws.Range("J15").Value = 0.630170773317832
$ws.Range("M15").Value = 0.8776213333333334
$ws.Range("N15").Value = 2.632864
$ws.Range("O15").Value = 0.126488678822602
$ws.Range("P15").Value = 0.126488678822602
$ws.Range("Q15").Value = 600.1662195984
$ws.Range("R15").Value = 5401.4959763856
$ws.Range("S15").Value = 0.07970946854958995
$ws.Range("T15").Value = 0.07970946854958995
$ws.Range("G16").Value = 683.8555499999999
$ws.Range("H16").Value = 2051.56665
$ws.Range("I16").Value = 0.630170773317832
$ws.Range("J16").Value = 0.630170773317832
$ws.Range("M16").Value = 5.010059000000001
$ws.Range("N16").Value = 15.030177
$ws.Range("O16").Value = 0.7220833401192995
$ws.Range("P16").Value = 0.7220833401192994
$ws.Range("Q16").Value = 3426.15665297745
$ws.Range("R16").Value = 30835.40987679705
$ws.Range("S16").Value = 0.4550358168429021
$ws.Range("T16").Value = 0.455035816842902
$ws.Range("G17").Value = 118.2996293333333
$ws.Range("H17").Value = 354.898888
$ws.Range("I17").Value = 0.1090127423842646
$ws.Range("J17").Value = 0.1090127423842646
$ws.Range("M17").Value = 1.050658666666667
$ws.Range("N17").Value = 3.151976
$ws.Range("O17").Value = 0.1514279810580986
$ws.Range("P17").Value = 0.1514279810580986
$ws.Range("Q17").Value = 124.2925308225209
$ws.Range("R17").Value = 1118.632777402688
$ws.Range("S17").Value = 0.0165075794888558
$ws.Range("T17").Value = 0.01650757948885579
$ws.Range("G18").Value = 118.2996293333333
$ws.Range("H18").Value = 354.898888
$ws.Range("I18").Value = 0.1090127423842646
$ws.Range("J18").Value = 0.1090127423842646
$ws.Range("M18").Value = 0.8776213333333334
$ws.Range("N18").Value = 2.632864
$ws.Range("O18").Value = 0.126488678822602
$ws.Range("P18").Value = 0.126488678822602
$ws.Range("Q18").Value = 103.8222784283591
$ws.Range("R18").Value = 934.400505855232
$ws.Range("S18").Value = 0.01378887775901429
$ws.Range("T18").Value = 0.01378887775901429
$ws.Range("G19").Value = 118.2996293333333
$ws.Range("H19").Value = 354.898888
$ws.Range("I19").Value = 0.1090127423842646
$ws.Range("J19").Value = 0.1090127423842646
$ws.Range("M19").Value = 5.010059000000001
$ws.Range("N19").Value = 15.030177
$ws.Range("O19").Value = 0.7220833401192995
$ws.Range("P19").Value = 0.7220833401192994
$ws.Range("Q19").Value = 592.6881226381307
$ws.Range("R19").Value = 5334.193103743177
$ws.Range("S19").Value = 0.07871628513639449
$ws.Range("T19").Value = 0.07871628513639448

Write-Output "Applied 216 cell updates"
